$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 75001
$ws.Range("E2").Value = "14625 FLANDERS CT"
$ws.Range("C3").Select()
